$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order Data")

# Add a new data row (row 11) to the Order Data sheet, mirroring the
# structure of the existing rows.
$ws.Cells.Item(11, 1).Value = 10.0
$ws.Cells.Item(11, 2).Value = "dineIn"
$ws.Cells.Item(11, 3).Value = "[1, 2, 2]"
$ws.Cells.Item(11, 4).Value = $false
$ws.Cells.Item(11, 5).Value = "InProgress"
$ws.Cells.Item(11, 6).Value = 4.0
$ws.Cells.Item(11, 7).Value = 0.0
$ws.Cells.Item(11, 8).Value = 0.0
$ws.Cells.Item(11, 9).Value = 0.0
